$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 1474
$ws.Cells.Item(3, 6).Value = 1445
$ws.Cells.Item(6, 6).Value = 712
$ws.Cells.Item(7, 6).Value = 36
$ws.Cells.Item(8, 6).Value = 642
$ws.Cells.Item(11, 6).Value = 1383
$ws.Cells.Item(12, 6).Value = 33605
$ws.Cells.Item(13, 6).Value = 7107
$ws.Cells.Item(14, 6).Value = 115
$ws.Cells.Item(15, 6).Value = 369
$ws.Cells.Item(16, 6).Value = 578
$ws.Cells.Item(17, 6).Value = 446
$ws.Cells.Item(19, 6).Value = 108
$ws.Cells.Item(20, 6).Value = 85
$ws.Cells.Item(22, 6).Value = 451
$ws.Cells.Item(23, 6).Value = 107
$ws.Cells.Item(24, 6).Value = 809
$ws.Cells.Item(25, 6).Value = 12
$ws.Cells.Item(26, 6).Value = 320
$ws.Cells.Item(28, 6).Value = 443
$ws.Cells.Item(29, 6).Value = 24
$ws.Cells.Item(30, 6).Value = 212
$ws.Cells.Item(31, 6).Value = 53
$ws.Cells.Item(32, 6).Value = 741
$ws.Cells.Item(33, 6).Value = 292
$ws.Cells.Item(35, 6).Value = 746
$ws.Cells.Item(36, 6).Value = 112
$ws.Cells.Item(38, 6).Value = 796
$ws.Cells.Item(39, 6).Value = 289

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 1209
$ws.Cells.Item(3, 6).Value = 8
$ws.Cells.Item(6, 6).Value = 291
$ws.Cells.Item(7, 6).Value = 4326
$ws.Cells.Item(9, 6).Value = 238
$ws.Cells.Item(12, 6).Value = 59
$ws.Cells.Item(16, 6).Value = 54

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 1463

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 1463
$ws.Cells.Item(4, 6).Value = 1209
$ws.Cells.Item(5, 6).Value = 1474
$ws.Cells.Item(6, 6).Value = 8
$ws.Cells.Item(7, 6).Value = 1445
$ws.Cells.Item(9, 6).Value = 712
$ws.Cells.Item(10, 6).Value = 36
$ws.Cells.Item(11, 6).Value = 642
$ws.Cells.Item(13, 6).Value = 1383
$ws.Cells.Item(15, 6).Value = 291
$ws.Cells.Item(16, 6).Value = 238
$ws.Cells.Item(17, 6).Value = 238
$ws.Cells.Item(20, 6).Value = 7107
$ws.Cells.Item(21, 6).Value = 369
$ws.Cells.Item(22, 6).Value = 59
$ws.Cells.Item(23, 6).Value = 578
$ws.Cells.Item(24, 6).Value = 446
$ws.Cells.Item(26, 6).Value = 108
$ws.Cells.Item(27, 6).Value = 86
$ws.Cells.Item(31, 6).Value = 451
$ws.Cells.Item(32, 6).Value = 107
$ws.Cells.Item(33, 6).Value = 809
$ws.Cells.Item(34, 6).Value = 12
$ws.Cells.Item(35, 6).Value = 320
$ws.Cells.Item(37, 6).Value = 444
$ws.Cells.Item(38, 6).Value = 24
$ws.Cells.Item(39, 6).Value = 212
$ws.Cells.Item(40, 6).Value = 53
$ws.Cells.Item(41, 6).Value = 741
$ws.Cells.Item(42, 6).Value = 54
$ws.Cells.Item(43, 6).Value = 292
$ws.Cells.Item(45, 6).Value = 796
$ws.Cells.Item(46, 6).Value = 289
